# "Generate Report for Archive"
#
# The localization-status report has moved on from the handoff stage:
# every "Ready for handoff" status cell becomes "In Translation" across
# the Overview sheet (per-language status columns) and each per-language
# detail sheet (zh-cn, de-de). Narrower status text means the Status
# column(s) no longer need to be as wide, so we also tighten those
# columns to their new auto-fit width.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newColumnWidth = 12.5   # renders to the narrower stored column width

# --- Overview sheet: zh-cn (E) and de-de (F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E1:E2").ColumnWidth = $newColumnWidth
$wsOverview.Range("F1:F2").ColumnWidth = $newColumnWidth

# --- zh-cn detail sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C1:C2").ColumnWidth = $newColumnWidth

# --- de-de detail sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C1:C2").ColumnWidth = $newColumnWidth
